$wb = $excel.ActiveWorkbook

# --- GeneralVariables sheet: append two new rows of data ---
$wsGeneral = $wb.Worksheets.Item("GeneralVariables")

$wsGeneral.Range("A15").Value = "nonMobileVoice"
$wsGeneral.Range("B15").Value = "Fix Voice Value Added Services"

$wsGeneral.Range("A16").Value = "optyStageCloseWon"
$wsGeneral.Range("B16").Value = "Closed Won"

# --- TC2 sheet: move its own selection, it's no longer the active tab ---
$wsTC2 = $wb.Worksheets.Item("TC2")
$wsTC2.Range("B2").Select() | Out-Null

# --- GeneralVariables becomes the active/selected tab ---
$wsGeneral.Activate() | Out-Null
$wsGeneral.Range("B17").Select() | Out-Null
